$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "64.446.44"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  +0.22%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.510.36"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  +0.04%  "

$ws.Range("E4").Value = "  +0.06%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "591.88"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +1.36%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "134.53"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -0.11%  "

$ws.Range("E7").Value = "  -0.01%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.488"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +0.12%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "7.54"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +6.37%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.124"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -0.04%  "

$ws.Range("E11").Value = "  +3.10%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "4.109.26"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +0.11%  "

$ws.Range("E13").Value = "  +1.60%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.0000181"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +0.86%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "3.510.72"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +0.06%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "25.79"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -3.40%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "64.422.95"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +0.20%  "

$ws.Range("E18").Value = "  +1.78%  "

$ws.Range("E19").Value = "  +3.07%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "13.63"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -1.51%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "394.34"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +2.69%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.575"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +1.29%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "3.651.69"
$ws.Range("D23").Style = "Normal"

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "74.74"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +1.05%  "

$ws.Range("E25").Value = "  +0.18%  "

$ws.Range("E26").Value = "  +0.34%  "

$ws.Range("E27").Value = "  +2.82%  "

$ws.Range("E28").Value = "  -0.03%  "

$ws.Range("E29").Value = "  -3.00%  "

$ws.Range("E30").Value = "  +2.00%  "

$ws.Range("E31").Value = "  -0.19%  "

$ws.Range("E32").Value = "  -5.91%  "

$ws.Range("E33").Value = "  +7.06%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "3.539.68"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +0.38%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "23.43"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -0.48%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "5.33"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +0.23%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "6.95"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +1.21%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "1.55"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +0.45%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "167.38"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +2.07%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.0791"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +0.81%  "

$ws.Range("E42").Value = "  +0.18%  "

$ws.Range("E43").Value = "  +0.07%  "

$ws.Range("B44").Value = "EnergySwap"
$ws.Range("C44").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "25.24"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -1.39%  "

$ws.Range("B45").Value = "Filecoin"
$ws.Range("C45").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "4.45"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +0.97%  "

$ws.Range("E46").Value = "  +2.99%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "1.17"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -3.08%  "

$ws.Range("E48").Value = "  +0.49%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "2.389.00"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -3.42%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.899"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -2.11%  "

$ws.Range("E51").Value = "  -0.43%  "
